$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: fix the casing of the "cost type" label -> "Cost type"
$ws.Range("D6").Value = "Cost type"

# Drop the two blank placeholder rows (7 and 8) under the table; that also
# retires the dedicated "light fill" style that only those rows used.
$ws.Rows("7:8").Delete() | Out-Null

# Give the four data columns explicit widths (previously default width).
$ws.Columns("A").ColumnWidth = 12.833333333333334
$ws.Columns("B").ColumnWidth = 21.5
$ws.Columns("C").ColumnWidth = 18.833333333333332
$ws.Columns("D").ColumnWidth = 22.333333333333332

# Move the saved cursor/selection to C9 (below the now-shorter table).
$ws.Range("C9").Select() | Out-Null
